# Refresh the "cryptos" price/volume table (Tue Apr 16 19:56:49 UTC 2024 run).
# Only columns D (Price) and E (Volume(1h)) change; both are plain text cells in the
# source sheet, so every write below targets Range.Value with a literal string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the new Price strings (e.g. "535.59", "1.00", "0.0790") are valid-looking
# numbers, so a plain .Value assignment would let Excel re-interpret them as numbers and
# silently drop formatting such as trailing zeros. Temporarily mark those cells as Text,
# write the literal string, then clear the temporary formatting again so the cells end up
# with the same (default/general) style as before the edit. Looping one address at a time
# (rather than building one multi-area Union Range) so every cell reliably gets the format.
$textCells = @('D5', 'D6', 'D11', 'D12', 'D14', 'D19', 'D20', 'D21', 'D22', 'D23', 'D24', 'D25', 'D26', 'D27', 'D29', 'D30', 'D32', 'D33', 'D35', 'D36', 'D37', 'D39', 'D40', 'D42', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '62.727.06'
$ws.Range('E2').Value = '  -1.19%  '
$ws.Range('D3').Value = '3.060.16'
$ws.Range('E3').Value = '  -1.48%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '535.59'
$ws.Range('E5').Value = '  -3.73%  '
$ws.Range('D6').Value = '132.23'
$ws.Range('E6').Value = '  -4.54%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '3.052.02'
$ws.Range('E8').Value = '  -1.48%  '
$ws.Range('E9').Value = '  -1.26%  '
$ws.Range('E10').Value = '  -4.22%  '
$ws.Range('D11').Value = '6.06'
$ws.Range('E11').Value = '  -10.05%  '
$ws.Range('D12').Value = '0.449'
$ws.Range('E12').Value = '  -1.53%  '
$ws.Range('E13').Value = '  +1.76%  '
$ws.Range('D14').Value = '34.02'
$ws.Range('E14').Value = '  -3.83%  '
$ws.Range('D15').Value = '3.552.37'
$ws.Range('E15').Value = '  -1.63%  '
$ws.Range('D16').Value = '62.744.34'
$ws.Range('E16').Value = '  -1.29%  '
$ws.Range('E17').Value = '  -0.63%  '
$ws.Range('D18').Value = '3.061.19'
$ws.Range('E18').Value = '  -1.40%  '
$ws.Range('D19').Value = '6.60'
$ws.Range('E19').Value = '  -1.82%  '
$ws.Range('D20').Value = '479.88'
$ws.Range('E20').Value = '  -5.50%  '
$ws.Range('D21').Value = '13.24'
$ws.Range('E21').Value = '  -3.42%  '
$ws.Range('D22').Value = '0.690'
$ws.Range('E22').Value = '  -3.18%  '
$ws.Range('D23').Value = '7.07'
$ws.Range('E23').Value = '  -4.38%  '
$ws.Range('D24').Value = '78.75'
$ws.Range('E24').Value = '  +0.82%  '
$ws.Range('D25').Value = '12.00'
$ws.Range('E25').Value = '  -3.73%  '
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.11%  '
$ws.Range('D27').Value = '2.69'
$ws.Range('E27').Value = '  -3.55%  '
$ws.Range('E28').Value = '  -2.95%  '
$ws.Range('D29').Value = '0.998'
$ws.Range('E29').Value = '  +0.13%  '
$ws.Range('D30').Value = '25.87'
$ws.Range('E30').Value = '  -1.86%  '
$ws.Range('E31').Value = '  -9.04%  '
$ws.Range('D32').Value = '1.10'
$ws.Range('E32').Value = '  -2.00%  '
$ws.Range('D33').Value = '56.83'
$ws.Range('E33').Value = '  -3.62%  '
$ws.Range('E34').Value = '  -7.61%  '
$ws.Range('D35').Value = '5.32'
$ws.Range('E35').Value = '  +2.43%  '
$ws.Range('D36').Value = '5.97'
$ws.Range('E36').Value = '  +0.60%  '
$ws.Range('D37').Value = '475.63'
$ws.Range('E37').Value = '  -12.01%  '
$ws.Range('D38').Value = '3.087.63'
$ws.Range('E38').Value = '  -0.26%  '
$ws.Range('D39').Value = '0.0391'
$ws.Range('E39').Value = '  -5.88%  '
$ws.Range('D40').Value = '0.0790'
$ws.Range('E40').Value = '  -1.61%  '
$ws.Range('E41').Value = '  -3.10%  '
$ws.Range('D42').Value = '8.04'
$ws.Range('E42').Value = '  -1.34%  '
$ws.Range('E43').Value = '  -2.11%  '
$ws.Range('E44').Value = '  -3.01%  '
$ws.Range('D46').Value = '0.0₃0532'
$ws.Range('E46').Value = '  +6.15%  '
$ws.Range('D47').Value = '120.98'
$ws.Range('E47').Value = '  -0.52%  '
$ws.Range('D48').Value = '2.00'
$ws.Range('E48').Value = '  -5.19%  '
$ws.Range('D49').Value = '24.24'
$ws.Range('E49').Value = '  +0.49%  '
$ws.Range('D50').Value = '0.108'
$ws.Range('E50').Value = '  +0.39%  '
$ws.Range('D51').Value = '2.28'
$ws.Range('E51').Value = '  -3.65%  '

foreach ($addr in $textCells) {
    $ws.Range($addr).ClearFormats()
}

